$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style from the last existing data row (A385) down
# through the new rows so column A keeps the same numFmt/border/alignment.
$ws.Range("A385").Copy()
$ws.Range("A386:A464").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New daily rows (date serial, nuovi pos., somma mobile 7gg., per 100k abitanti)
$data = @"
386,44460,4,31,121.3354730126424
387,44461,4,35,136.9916630787898
388,44462,4,30,117.4214254961055
389,44463,5,34,133.0776155622529
390,44464,1,29,113.5073779795687
391,44465,4,24,93.93714039688442
392,44466,2,24,93.93714039688442
393,44467,4,24,93.93714039688442
394,44468,0,20,78.28095033073701
395,44469,4,20,78.28095033073701
396,44470,3,18,70.45285529766332
397,44471,0,17,66.53880778112647
398,44472,2,15,58.71071274805276
399,44473,2,15,58.71071274805276
400,44474,4,15,58.71071274805276
401,44475,1,16,62.62476026458961
402,44476,2,14,54.79666523151591
403,44477,1,12,46.96857019844221
404,44478,0,12,46.96857019844221
405,44479,0,10,39.14047516536851
406,44480,2,10,39.14047516536851
407,44481,0,6,23.4842850992211
408,44482,0,5,19.57023758268425
409,44483,0,3,11.74214254961055
410,44484,0,2,7.828095033073701
411,44485,0,2,7.828095033073701
412,44486,0,2,7.828095033073701
413,44487,4,4,15.6561900661474
414,44488,0,4,15.6561900661474
415,44489,1,5,19.57023758268425
416,44490,4,9,35.22642764883166
417,44491,0,9,35.22642764883166
418,44492,2,11,43.05452268190535
419,44493,0,11,43.05452268190535
420,44494,2,9,35.22642764883166
421,44495,1,10,39.14047516536851
422,44496,0,9,35.22642764883166
423,44497,0,5,19.57023758268425
424,44498,1,6,23.4842850992211
425,44499,1,5,19.57023758268425
426,44500,2,7,27.39833261575795
427,44501,0,5,19.57023758268425
428,44502,1,5,19.57023758268425
429,44503,0,5,19.57023758268425
430,44504,0,5,19.57023758268425
431,44505,5,9,35.22642764883166
432,44506,4,12,46.96857019844221
433,44507,0,10,39.14047516536851
434,44508,0,10,39.14047516536851
435,44509,0,9,35.22642764883166
436,44510,1,10,39.14047516536851
437,44511,1,11,43.05452268190535
438,44512,2,8,31.31238013229481
439,44513,1,5,19.57023758268425
440,44514,1,6,23.4842850992211
441,44515,5,11,43.05452268190535
442,44516,30,41,160.4759481780109
443,44517,2,42,164.3899956945477
444,44518,2,43,168.3040432110846
445,44519,3,44,172.2180907276214
446,44520,6,49,191.7883283103057
447,44521,2,50,195.7023758268425
448,44522,9,54,211.3585658929899
449,44523,11,35,136.9916630787898
450,44524,14,47,183.960233277232
451,44525,10,55,215.2726134095268
452,44526,9,61,238.7568985087479
453,44527,5,60,234.842850992211
454,44528,7,65,254.4130885748953
455,44529,3,59,230.9288034756742
456,44530,14,62,242.6709460252848
457,44531,0,48,187.8742807937688
458,44532,8,46,180.0461857606951
459,44533,11,48,187.8742807937688
460,44534,8,51,199.6164233433794
461,44535,18,62,242.6709460252848
462,44536,15,74,289.639516223727
463,44537,10,70,273.9833261575795
464,44538,11,81,317.0378488394849
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 1).Value = [double]$parts[1]
    $ws.Cells.Item($r, 2).Value = [double]$parts[2]
    $ws.Cells.Item($r, 3).Value = [double]$parts[3]
    $ws.Cells.Item($r, 4).Value = [double]$parts[4]
}
